$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values per repull of data / mean calculation fix
$ws.Range("F2").Value = -5
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -9
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = 4
$ws.Range("F19").Value = 2
$ws.Range("F22").Value = -2
$ws.Range("F24").Value = -1
